# "Move all constants to PhysicalConstants"
#
# A new parameter row ("rpmax" / "Max r prime") is inserted into the
# Source block (right after the existing "SourceMode" row, i.e. before
# the old row 6), pushing every following row down by one.
#
# Old row 6 ("Capture"/"Drift"/Length) becomes new row 7, and so on down
# to old row 23 becoming new row 24. Excel's native row-insert takes care
# of shifting all the existing cell content/formulas/styles down and
# re-pointing relative formula references (F10, F11, F19, F20, I6, I8,
# I12, I13, I15, I17, I18, I22 -> F11, F12, F20, F21, I7, I9, I13, I14,
# I16, I18, I19, I23) automatically, exactly like it does in the real
# Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, empty row at row 6 - everything at/after row 6 shifts
# down by one (row 6 -> 7, row 23 -> 24, dimension A1:I23 -> A1:I24).
$ws.Rows.Item(6).Insert()

# Row 15 (formerly row 14, the last "Capture"/"Aperture" row, which
# closes out a sub-block with a full thin border on every side) has the
# exact formatting the new row needs. Copy its formats only - not its
# values/formulas - onto the freshly inserted row 6.
$ws.Range("A15:H15").Copy()
$ws.Range("A6:H6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row: Stage=1, Section/Element=Source/Source (same
# block as row 5), Type=Parameterised TNSA, new Parameter "rpmax",
# Value 0.005, no Unit, Comment "Max r prime".
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Source"
$ws.Range("C6").Value = "Source"
$ws.Range("D6").Value = "Parameterised TNSA"
$ws.Range("E6").Value = "rpmax"
$ws.Range("F6").Value = 0.005
$ws.Range("H6").Value = "Max r prime"
